# Project 3 Plan.xlsx — "Corrected fileClean.py and initial plotting"
#
# The upstream data-cleaning script dropped the standalone "Flask Server (API)"
# task and re-sequenced / re-percented several of the remaining tasks. This
# script reproduces that on the live workbook:
#   1. Removes the "Flask Server (API)" activity (by overwriting rows 12-18
#      with the corrected activity/date/owner/% data, which leaves that
#      shared string unreferenced so Excel drops it on save) and clears the
#      now-unused last row (19).
#   2. Updates the "% Complete" values that changed.
#   3. Re-stripes the Gantt bar helper cells (columns G:Q) so the shaded
#      "elapsed" cells line up with the new start/duration/percent-complete
#      numbers, by copying the format of a cell that already carries the
#      desired look and pasting formats-only onto the cells that need it.
#   4. Leaves the selection on B2, matching the saved view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update activity/date/percent/owner values for the corrected schedule ---
$ws.Range("E8").Value = 1
$ws.Range("E11").Value = 1

$ws.Range("B12").Value = "Code Reviews"
$ws.Range("C12").Value = 45393
$ws.Range("D12").Value = 45399
$ws.Range("E12").Value = 0.1
$ws.Range("F12").Value = "Brennan"

$ws.Range("B13").Value = "Visualization Page"
$ws.Range("C13").Value = 45396
$ws.Range("D13").Value = 45398
$ws.Range("E13").Value = 0.1
$ws.Range("F13").Value = "Tim"

$ws.Range("B14").Value = "Presentation Template"
$ws.Range("C14").Value = 45394
$ws.Range("D14").Value = 45397
$ws.Range("F14").Value = "Kortney"

$ws.Range("B15").Value = "Visualization QA"
$ws.Range("C15").Value = 45398
$ws.Range("D15").Value = 45399

$ws.Range("B16").Value = "Project ReadMe"
$ws.Range("C16").Value = 45399
$ws.Range("E16").Value = 0.1
$ws.Range("F16").Value = "Tim"

$ws.Range("B17").Value = "Project Review"
$ws.Range("F17").Value = "Team"

$ws.Range("B18").Value = "Project Presentation"
$ws.Range("C18").Value = 45400
$ws.Range("D18").Value = 45400

# Row 19 ("Project Presentation") no longer exists as a separate activity;
# its old data moved up into row 18, so blank the row out entirely.
$ws.Range("B19:Q19").Clear()

# --- Re-stripe the Gantt "status" cells (fill style indices) to match the new schedule ---
$ws.Range("G5").Copy()
$ws.Range("H8:I8").PasteSpecial(-4122)
$ws.Range("I11:K11").PasteSpecial(-4122)
$ws.Range("I12").PasteSpecial(-4122)
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("I16").PasteSpecial(-4122)

$ws.Range("H5").Copy()
$ws.Range("P14").PasteSpecial(-4122)
$ws.Range("L15").PasteSpecial(-4122)
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("P18").PasteSpecial(-4122)

$ws.Range("I5").Copy()
$ws.Range("O14").PasteSpecial(-4122)
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("M15").PasteSpecial(-4122)

$ws.Range("J12").Copy()
$ws.Range("M12:P12").PasteSpecial(-4122)
$ws.Range("J14:L14").PasteSpecial(-4122)
$ws.Range("O15:P15").PasteSpecial(-4122)
$ws.Range("J16:N16").PasteSpecial(-4122)
$ws.Range("Q18").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Restore the saved selection ---
$ws.Range("B2").Select()
